$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.088.88'
$ws.Range("E2").Value = '  -0.08%  '
$ws.Range("D3").Value = '1.624.70'
$ws.Range("E3").Value = '  -0.84%  '
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.34'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.99%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.518'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.47%  '
$ws.Range("E7").Value = '  -0.13%  '
$ws.Range("E8").Value = '  +0.81%  '
$ws.Range("E9").Value = '  -1.34%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.04'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.30%  '
$ws.Range("E11").Value = '  -0.08%  '
$ws.Range("D12").Value = '1.851.54'
$ws.Range("E12").Value = '  -0.90%  '
$ws.Range("D13").Value = '1.622.22'
$ws.Range("E13").Value = '  -1.97%  '
$ws.Range("E14").Value = '  +0.20%  '
$ws.Range("E15").Value = '  +0.15%  '
$ws.Range("B16").Value = 'Litecoin'
$ws.Range("C16").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '64.56'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.11%  '
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = '27.062.25'
$ws.Range("E17").Value = '  -0.27%  '
$ws.Range("D18").Value = '0.0₃0739'
$ws.Range("E18").Value = '  +0.09%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '214.23'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.20%  '
$ws.Range("E20").Value = '  -0.05%  '
$ws.Range("E21").Value = '  -1.30%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.35'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.24%  '
$ws.Range("E23").Value = '  -6.66%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.05'
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '148.35'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.11%  '
$ws.Range("E26").Value = '  -0.12%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.38'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.38%  '
$ws.Range("E28").Value = '  -2.80%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.55'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.78%  '
$ws.Range("E30").Value = '  +0.41%  '
$ws.Range("E31").Value = '  -0.70%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.35'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.85%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.731'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +34.89%  '
$ws.Range("E34").Value = '  -0.36%  '
$ws.Range("D35").Value = '1.358.53'
$ws.Range("E35").Value = '  +4.02%  '
$ws.Range("E36").Value = '  -0.19%  '
$ws.Range("E37").Value = '  -0.68%  '
$ws.Range("E38").Value = '  +1.00%  '
$ws.Range("E39").Value = '  -1.31%  '
$ws.Range("E40").Value = '  -0.16%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.802'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.20%  '
$ws.Range("E42").Value = '  +0.43%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '64.27'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.96%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.33'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.74%  '
$ws.Range("D45").Value = '1.762.74'
$ws.Range("E45").Value = '  -0.90%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.66'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.74%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '89.98'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.51%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.866'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +29.79%  '
$ws.Range("D49").Value = '0.0₆0105'
$ws.Range("E49").Value = '  -2.09%  '
$ws.Range("E50").Value = '  +5.04%  '
$ws.Range("E51").Value = '  +0.12%  '
